# [DataTable] FarmerStatTable -> FarmerLevelTable 마이그레이션 & 테이블 정리
#
# Insert a new "type" row (id/level/costItemID/costValue all -> "int")
# right after the header row, pushing the existing data rows down by one,
# then restore the freeze-panes/selection view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the first data row (row 3) - this shifts all the
# existing data rows (old rows 3-57) down to rows 4-58.
$ws.Rows.Item(3).Insert()

# Fill the freshly inserted row with the shared "int" label in every column.
$ws.Range("A3:D3").Value = "int"

# Copy the header row's formatting (fill/border/font/alignment) onto the new
# row so it reads as a second header/type row.
$ws.Range("A2:D2").Copy()
$ws.Range("A3:D3").PasteSpecial(-4122)  # xlPasteFormats

# Both the header row and the new type row grew a bit taller.
$ws.Rows.Item(2).RowHeight = 19
$ws.Rows.Item(3).RowHeight = 19

# Re-establish freeze panes at C4 (2 columns / 3 rows frozen) with the same
# per-pane selections as the saved view.
$null = $ws.Range("C4").Select()
$excel.ActiveWindow.FreezePanes = $true
$null = $ws.Range("P8").Select()
